# Weekly update: prepend a new observation row for Ciboulette / Vega Modelo de Temuco.
# Insert a new row above row 60, shifting existing rows 60:178 down to 61:179.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("60:60").Insert(-4121)

$newRow = $ws.Rows("60:60")
$newRow.Cells.Item(1, 1).Value2 = 10
$newRow.Cells.Item(1, 2).Value2 = "Vega Modelo de Temuco"
$newRow.Cells.Item(1, 3).Value2 = "La Araucanía"
$newRow.Cells.Item(1, 4).Value2 = 44519
$newRow.Cells.Item(1, 5).Value2 = 9
$newRow.Cells.Item(1, 6).Value2 = 100112039
$newRow.Cells.Item(1, 7).Value2 = "Ciboulette"
$newRow.Cells.Item(1, 8).Value2 = "Sin especificar"
$newRow.Cells.Item(1, 9).Value2 = "Primera"
$newRow.Cells.Item(1, 10).Value2 = 65
$newRow.Cells.Item(1, 11).Value2 = 5000
$newRow.Cells.Item(1, 12).Value2 = 5000
$newRow.Cells.Item(1, 13).Value2 = 5000
$newRow.Cells.Item(1, 14).Value2 = "`$/docena de atados"
$newRow.Cells.Item(1, 15).Value2 = "Provincia de Cautín"
$newRow.Cells.Item(1, 16).Value2 = 1667
$newRow.Cells.Item(1, 17).Value2 = 3
$newRow.Cells.Item(1, 18).Value2 = "Hortaliza"
